$wb = $excel.ActiveWorkbook

# --- "login" sheet: rows 2..21, columns G,H (username) and I (email) ---
# Usernames go from "<Name>18" -> "<Name>19"
$loginSheet = $wb.Worksheets.Item("login")
for ($r = 2; $r -le 21; $r++) {
    $g = $loginSheet.Cells.Item($r, 7).Value2   # column G
    $i = $loginSheet.Cells.Item($r, 9).Value2   # column I
    if ($g -ne $null -and $g.ToString().EndsWith("18")) {
        $newName = $g.ToString().Substring(0, $g.ToString().Length - 2) + "19"
        $loginSheet.Cells.Item($r, 7).Value2 = $newName
        $loginSheet.Cells.Item($r, 8).Value2 = $newName
    }
    if ($i -ne $null -and $i.ToString().Contains("18@gmail.com")) {
        $newEmail = $i.ToString().Replace("18@gmail.com", "19@gmail.com")
        $loginSheet.Cells.Item($r, 9).Value2 = $newEmail
    }
}

# --- "order" sheet: rows 2..21, columns R,S (username) and T (email) ---
# Usernames go from "<Name>23" -> "<Name>24"
$orderSheet = $wb.Worksheets.Item("order")
for ($r = 2; $r -le 21; $r++) {
    $rCol = $orderSheet.Cells.Item($r, 18).Value2  # column R
    $tCol = $orderSheet.Cells.Item($r, 20).Value2  # column T
    if ($rCol -ne $null -and $rCol.ToString().EndsWith("23")) {
        $newName = $rCol.ToString().Substring(0, $rCol.ToString().Length - 2) + "24"
        $orderSheet.Cells.Item($r, 18).Value2 = $newName
        $orderSheet.Cells.Item($r, 19).Value2 = $newName
    }
    if ($tCol -ne $null -and $tCol.ToString().Contains("23@gmail.com")) {
        $newEmail = $tCol.ToString().Replace("23@gmail.com", "24@gmail.com")
        $orderSheet.Cells.Item($r, 20).Value2 = $newEmail
    }
}

# --- "Sheet1" sheet: the two driver cells I2 and I23 feed CONCATENATE formulas ---
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Cells.Item(2, 9).Value2 = 19
$sheet1.Cells.Item(23, 9).Value2 = 24

# Scroll position change on Sheet1's sheet view
$sheet1.Activate()
$excel.Goto($sheet1.Range("A18"), $true)
